# Add two new worksheets ("Скидка за представленность" and
# "Скидка на РБ от закупа продукци") at the end of the workbook, mirroring
# the layout already used by the "Скидка за объем закупа" sheet (headers +
# one data row + a totals row), and make the last one the active sheet.

$wb = $excel.ActiveWorkbook

# Reference sheet whose header/total formatting (fill + border, style index
# 1) we want to reuse for the new sheets' header and totals rows.
$wsRef = $wb.Worksheets.Item(2)

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

# --- New sheet 4: "Скидка за представленность" ----------------------------
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Скидка за представленность"

$ws4.Range("A1").Value = "Период"
$ws4.Range("B1").Value = "Номер договора/ДС"
$ws4.Range("C1").Value = "Тип скидки"
$ws4.Range("D1").Value = "Скидка %"
$ws4.Range("E1").Value = "Сумма скидки"

$ws4.Range("A2").Value = "01.03.2022-01.01.2022"
$ws4.Range("B2").Value = "2500800DLR"
$ws4.Range("C2").Value = "Скидка за представленность"
$ws4.Range("D2").Value = 12
$ws4.Range("E2").Value = 19696.32

$ws4.Range("D3").Value = "Итог:"
$ws4.Range("E3").Value = 19696.32

$wsRef.Range("A1").Copy()
$ws4.Range("A1:E1").PasteSpecial(-4122)
$ws4.Range("D3:E3").PasteSpecial(-4122)

# --- New sheet 5: "Скидка на РБ от закупа продукци" ------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet2)
$ws5.Name = "Скидка на РБ от закупа продукци"

$ws5.Range("A1").Value = "Период"
$ws5.Range("B1").Value = "Номер договора/ДС"
$ws5.Range("C1").Value = "Тип скидки"
$ws5.Range("D1").Value = "Скидка %"
$ws5.Range("E1").Value = "Сумма скидки"

$ws5.Range("A2").Value = "01.03.2022-01.01.2022"
$ws5.Range("B2").Value = "2500800DLR"
$ws5.Range("C2").Value = "Скидка за прСкидка на РБ от закупа продукциедставленность"
$ws5.Range("D2").Value = 12
$ws5.Range("E2").Value = 19696.32

$ws5.Range("D3").Value = "Итог:"
$ws5.Range("E3").Value = 19696.32

$wsRef.Range("A1").Copy()
$ws5.Range("A1:E1").PasteSpecial(-4122)
$ws5.Range("D3:E3").PasteSpecial(-4122)

# Restore the per-sheet cursor/selection state and make the newly added
# last sheet the active / selected tab.
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

[void]$ws1.Range("E25").Select()
[void]$ws2.Range("A1").Select()
[void]$ws3.Range("A1:G5").Select()
[void]$ws4.Range("A1:E3").Select()
[void]$ws5.Range("E10").Select()

[void]$ws5.Activate()
